# Append new scraped rows to the top of the list (most-recent-first ordering),
# refresh the "taken at" timestamp on every row, and add one more row at the
# bottom for a listing that previously had no skill-match tag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-01-07 12:41:04"

# ---------------------------------------------------------------------------
# 1) Make room: insert 4 fresh rows above the existing data (old rows 2-4
#    become rows 6-8). Hyperlinks do not follow a row insert automatically,
#    so we rebuild the hyperlinks collection from scratch at the end.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Range("A2:A5").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2) Data for every data row (2-9) after the insert, in final position.
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row = 2;  A = $newTimestamp; B = "産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5450864"; G = 383; H = "🔥AI,Ai ◆開発" },
    @{ Row = 3;  A = $newTimestamp; B = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"; C = "システム開発"; D = "5,000 円 ~ 10,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5466794"; G = 135; H = "◆ツール,スクレイピング ◇サイト" },
    @{ Row = 4;  A = $newTimestamp; B = "【大募集】iOSアプリ開発!Apple開発者アカウントも一緒に譲渡希望です。"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5466844"; G = 88; H = "◆開発 ◇アプリ" },
    @{ Row = 5;  A = $newTimestamp; B = "FileMaker開発"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5466845"; G = 68; H = "◆開発" },
    @{ Row = 6;  A = $newTimestamp; B = "Raspberry Piでの開発"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5466611"; G = 68; H = "◆開発" },
    @{ Row = 7;  A = $newTimestamp; B = "【急募】GBP一括投稿システムのインスタ連携改修依頼"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5466476"; G = 33; H = $null },
    @{ Row = 8;  A = $newTimestamp; B = "【長期/業務委託】UX改善をリードできるフロント寄り Laravel エンジニア募集(リモート可)"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5466459"; G = 25; H = $null },
    @{ Row = 9;  A = $newTimestamp; B = "【EC-CUBE】定期購入機能の調査・改修依頼"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5466925"; G = 18; H = $null }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    if ($r.H -ne $null) {
        $ws.Cells.Item($row, 8).Value = $r.H
    }

    # Recreate the hyperlink on the URL cell and restore the Hyperlink style.
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $r.F)
    $ws.Cells.Item($row, 6).Style = "Hyperlink"
}

# ---------------------------------------------------------------------------
# 3) Widen column H slightly to fit the new, longer skill-tag text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 18.166666666666668
